$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to Text
# format first, otherwise Excel will auto-convert them into numeric values
# instead of preserving them as text (matching the source report formatting).
$numericLookingRefs = @("D5","D6","D8","D9","D10","D11","D14","D15","D16","D19","D21","D22","D25","D26","D27","D28","D29","D33","D35","D36","D37","D38","D40","D41","D42","D43","D44","D46","D49","D51")
foreach ($ref in $numericLookingRefs) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply the updated values cell by cell
$ws.Range("D2").Value = "34.961.47"
$ws.Range("E2").Value = "  -1.66%  "
$ws.Range("D3").Value = "1.890.78"
$ws.Range("E3").Value = "  -1.06%  "
$ws.Range("E4").Value = "  -0.29%  "
$ws.Range("D5").Value = "250.08"
$ws.Range("D6").Value = "0.691"
$ws.Range("E6").Value = "  -1.88%  "
$ws.Range("E7").Value = "  -0.26%  "
$ws.Range("D8").Value = "41.18"
$ws.Range("E8").Value = "  +0.36%  "
$ws.Range("D9").Value = "0.350"
$ws.Range("E9").Value = "  -1.77%  "
$ws.Range("D10").Value = "51.15"
$ws.Range("E10").Value = "  -2.72%  "
$ws.Range("D11").Value = "0.0741"
$ws.Range("E11").Value = "  +1.10%  "
$ws.Range("E12").Value = "  -1.68%  "
$ws.Range("D13").Value = "2.165.56"
$ws.Range("E13").Value = "  -0.98%  "
$ws.Range("D14").Value = "12.87"
$ws.Range("E14").Value = "  +2.13%  "
$ws.Range("D15").Value = "0.718"
$ws.Range("E15").Value = "  +0.18%  "
$ws.Range("D16").Value = "4.92"
$ws.Range("E16").Value = "  -0.23%  "
$ws.Range("D17").Value = "1.907.57"
$ws.Range("E17").Value = "  -0.19%  "
$ws.Range("D18").Value = "34.974.28"
$ws.Range("E18").Value = "  -1.52%  "
$ws.Range("D19").Value = "73.35"
$ws.Range("E19").Value = "  -0.08%  "
$ws.Range("D20").Value = "0.0₃0826"
$ws.Range("E20").Value = "  -0.39%  "
$ws.Range("D21").Value = "248.78"
$ws.Range("E21").Value = "  +2.35%  "
$ws.Range("D22").Value = "12.80"
$ws.Range("E22").Value = "  -3.18%  "
$ws.Range("E23").Value = "  -2.36%  "
$ws.Range("E24").Value = "  -0.28%  "
$ws.Range("D25").Value = "2.40"
$ws.Range("E25").Value = "  +3.59%  "
$ws.Range("D26").Value = "2.22"
$ws.Range("E26").Value = "  -4.37%  "
$ws.Range("D27").Value = "165.49"
$ws.Range("E27").Value = "  -2.60%  "
$ws.Range("D28").Value = "8.44"
$ws.Range("E28").Value = "  -3.12%  "
$ws.Range("D29").Value = "18.29"
$ws.Range("E29").Value = "  -3.10%  "
$ws.Range("E30").Value = "  -3.25%  "
$ws.Range("D31").Value = "4.128.64"
$ws.Range("E31").Value = "  +0.47%  "
$ws.Range("E32").Value = "  +0.75%  "
$ws.Range("D33").Value = "0.0584"
$ws.Range("E33").Value = "  +1.40%  "
$ws.Range("E34").Value = "  +5.30%  "
$ws.Range("B35").Value = "InternetComputer(DFINITY)"
$ws.Range("C35").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D35").Value = "4.18"
$ws.Range("E35").Value = "  -0.96%  "
$ws.Range("B36").Value = "BinanceUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D36").Value = "1.01"
$ws.Range("E36").Value = "  -0.35%  "
$ws.Range("B37").Value = "WEMIXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D37").Value = "1.85"
$ws.Range("E37").Value = "  -1.02%  "
$ws.Range("D38").Value = "0.840"
$ws.Range("E38").Value = "  -8.31%  "
$ws.Range("E39").Value = "  -1.72%  "
$ws.Range("D40").Value = "17.41"
$ws.Range("E40").Value = "  +0.86%  "
$ws.Range("D41").Value = "98.31"
$ws.Range("E41").Value = "  +0.30%  "
$ws.Range("D42").Value = "0.0669"
$ws.Range("E42").Value = "  +2.72%  "
$ws.Range("D43").Value = "0.0211"
$ws.Range("E43").Value = "  +0.65%  "
$ws.Range("D44").Value = "1.09"
$ws.Range("E44").Value = "  -3.09%  "
$ws.Range("D45").Value = "1.296.72"
$ws.Range("E45").Value = "  -4.60%  "
$ws.Range("D46").Value = "2.37"
$ws.Range("E46").Value = "  -2.66%  "
$ws.Range("E47").Value = "  -0.09%  "
$ws.Range("E48").Value = "  -1.73%  "
$ws.Range("D49").Value = "12.24"
$ws.Range("E49").Value = "  +0.06%  "
$ws.Range("E50").Value = "  +6.33%  "
$ws.Range("D51").Value = "6.48"
$ws.Range("E51").Value = "  -1.14%  "
